# Add Few More TestCases And Refactor Structure

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("multiSelectSelectAllValueData")
$ws3 = $wb.Worksheets.Item("multiSelectSelectFirstValueData")

# Work first on sheet2 so its selection stays A18 but it is no longer the active tab
$ws2.Activate()
$ws2.Range("A18").Select()

# Work on the third sheet: multiSelectSelectFirstValueData
$ws3.Activate()

# Add the new test case value in the next empty row (A3)
$ws3.Range("A3").Value = "Failed for Raport Purpoose"

# Update selection on sheet3 to the newly added cell and make it the active tab
$ws3.Range("A3").Select()
